$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update PoS pretrain condition result values
$ws.Range("C2").Value = 0.9787187739463602
$ws.Range("D2").Value = 0.810360153256705
$ws.Range("E2").Value = 0.7337164750957854
$ws.Range("F2").Value = 0.8548934865900383
$ws.Range("G2").Value = 0.8065409961685823
$ws.Range("H2").Value = 0.6467310344827586
$ws.Range("I2").Value = 0.557391570881226
$ws.Range("K2").Value = 0.7021180076628353
$ws.Range("L2").Value = 0.6832613026819924
$ws.Range("M2").Value = 0.6263785440613027
$ws.Range("N2").Value = 0.5810298850574712
$ws.Range("O2").Value = 0.5382283524904214
$ws.Range("P2").Value = 0.6748842911877395
$ws.Range("Q2").Value = 0.6471846743295019
$ws.Range("R2").Value = 0.7383662835249042
$ws.Range("C3").Value = 0.8256851196000076
$ws.Range("D3").Value = 0.9869131266453287
$ws.Range("E3").Value = 0.8633359216681502
$ws.Range("F3").Value = 0.8864983617734513
$ws.Range("G3").Value = 0.7844169617999659
$ws.Range("H3").Value = 0.5815041381792011
$ws.Range("I3").Value = 0.5821480653775497
$ws.Range("K3").Value = 0.8327115016761046
$ws.Range("L3").Value = 0.626048749076722
$ws.Range("M3").Value = 0.5896668623700309
$ws.Range("N3").Value = 0.5204068104770743
$ws.Range("O3").Value = 0.4949527471070624
$ws.Range("P3").Value = 0.6242684797636409
$ws.Range("Q3").Value = 0.6516353856934528
$ws.Range("R3").Value = 0.7544175299710233
$ws.Range("K4").Value = 0.7619742093951489
$ws.Range("K5").Value = 0.7897283588066475
$ws.Range("K6").Value = 0.6918146859549533
$ws.Range("C7").Value = 0.5348817848817848
$ws.Range("D7").Value = 0.4408924408924409
$ws.Range("E7").Value = 0.3677156177156177
$ws.Range("F7").Value = 0.6205461205461206
$ws.Range("G7").Value = 0.6103896103896104
$ws.Range("H7").Value = 0.9638694638694638
$ws.Range("I7").Value = 0.5282217782217782
$ws.Range("K7").Value = 0.571012321012321
$ws.Range("L7").Value = 0.5994838494838495
$ws.Range("M7").Value = 0.5815850815850816
$ws.Range("N7").Value = 0.6424408924408924
$ws.Range("O7").Value = 0.6506826506826506
$ws.Range("P7").Value = 0.5231435231435232
$ws.Range("Q7").Value = 0.5092407592407593
$ws.Range("R7").Value = 0.6644189144189144
$ws.Range("K8").Value = 0.6311166875784191
$ws.Range("K9").Value = 0.6224800645103485
$ws.Range("C10").Value = 0.6864176570458405
$ws.Range("D10").Value = 0.8378607809847198
$ws.Range("E10").Value = 0.7747877758913413
$ws.Range("F10").Value = 0.8368421052631579
$ws.Range("G10").Value = 0.6973684210526315
$ws.Range("H10").Value = 0.6874363327674023
$ws.Range("I10").Value = 0.5953310696095077
$ws.Range("K10").Value = 0.9286078098471986
$ws.Range("L10").Value = 0.749660441426146
$ws.Range("M10").Value = 0.6971986417657046
$ws.Range("N10").Value = 0.6166383701188455
$ws.Range("O10").Value = 0.582258064516129
$ws.Range("P10").Value = 0.6729202037351443
$ws.Range("Q10").Value = 0.6131578947368421
$ws.Range("R10").Value = 0.7716468590831919
$ws.Range("C11").Value = 0.77695730379627
$ws.Range("D11").Value = 0.851415317618101
$ws.Range("E11").Value = 0.7652182145223895
$ws.Range("F11").Value = 0.8307299062766259
$ws.Range("G11").Value = 0.8305879011644419
$ws.Range("H11").Value = 0.7406986651519455
$ws.Range("I11").Value = 0.5973681719208558
$ws.Range("K11").Value = 0.7626147874656821
$ws.Range("L11").Value = 0.9675754993846445
$ws.Range("M11").Value = 0.8019975385780554
$ws.Range("N11").Value = 0.6774590551926536
$ws.Range("O11").Value = 0.7219066553062576
$ws.Range("P11").Value = 0.7824008330966582
$ws.Range("Q11").Value = 0.6713055003313453
$ws.Range("R11").Value = 0.8437943765975575
$ws.Range("K12").Value = 0.6957003364240584
$ws.Range("K13").Value = 0.5687090178382571
$ws.Range("C14").Value = 0.3985592357685381
$ws.Range("D14").Value = 0.3056142823584684
$ws.Range("E14").Value = 0.2662281732049174
$ws.Range("F14").Value = 0.4128885756792733
$ws.Range("G14").Value = 0.4604964372406233
$ws.Range("H14").Value = 0.5389554459321901
$ws.Range("I14").Value = 0.3994205622112599
$ws.Range("K14").Value = 0.3534570511314697
$ws.Range("L14").Value = 0.4625322997416021
$ws.Range("M14").Value = 0.5369978858350951
$ws.Range("N14").Value = 0.58656330749354
$ws.Range("O14").Value = 0.9791715605669093
$ws.Range("P14").Value = 0.5359799545846058
$ws.Range("Q14").Value = 0.3555712160363323
$ws.Range("R14").Value = 0.4778795709028267
$ws.Range("K15").Value = 0.6893895065340355
$ws.Range("K16").Value = 0.6246587956720492
$ws.Range("K17").Value = 0.6235713813833652
$ws.Range("C19").Value = 0.6511833641613349
$ws.Range("D19").Value = 0.6875341039700426
$ws.Range("E19").Value = 0.6394785991363253
$ws.Range("F19").Value = 0.6936976124980982
$ws.Range("G19").Value = 0.6864967638709107
$ws.Range("H19").Value = 0.6234215045607119
$ws.Range("I19").Value = 0.5526112959795796
$ws.Range("K19").Value = 0.6614037808697132
$ws.Range("L19").Value = 0.6402403987741813
$ws.Range("M19").Value = 0.6217551209873585
$ws.Range("N19").Value = 0.5721452488883779
$ws.Range("O19").Value = 0.572062000368799
$ws.Range("P19").Value = 0.613324325339295
$ws.Range("Q19").Value = 0.6058508911559701
$ws.Range("R19").Value = 0.7127692949960444
$ws.Range("C24").Value = 0.8355683181540132
$ws.Range("D24").Value = 0.843948512539462
$ws.Range("E24").Value = 0.7900824301343616
$ws.Range("F24").Value = 0.8442304659166722
$ws.Range("G24").Value = 0.8036671029611815
$ws.Range("H24").Value = 0.6350996109242726
$ws.Range("I24").Value = 0.5649516931878795
$ws.Range("K24").Value = 0.7556693526991379
$ws.Range("L24").Value = 0.7031024411549016
$ws.Range("M24").Value = 0.6422269487762176
$ws.Range("N24").Value = 0.5694710925400618
$ws.Range("O24").Value = 0.5369960812746044
$ws.Range("P24").Value = 0.6686089213326457
$ws.Range("Q24").Value = 0.6857143806573156
$ws.Range("R24").Value = 0.7533063541494186
$ws.Range("C25").Value = 0.5614852180423321
$ws.Range("D25").Value = 0.6113602155316298
$ws.Range("E25").Value = 0.5759914634401546
$ws.Range("F25").Value = 0.6665896165919227
$ws.Range("G25").Value = 0.6272389856644431
$ws.Range("H25").Value = 0.6492462131335274
$ws.Range("I25").Value = 0.5659729813530742
$ws.Range("K25").Value = 0.6082030243670296
$ws.Range("L25").Value = 0.6320186456832088
$ws.Range("M25").Value = 0.6095290100485518
$ws.Range("N25").Value = 0.574460446489103
$ws.Range("O25").Value = 0.5758815360533673
$ws.Range("P25").Value = 0.5573151628394422
$ws.Range("Q25").Value = 0.5764490937807709
$ws.Range("R25").Value = 0.6926865793957868
$ws.Range("C26").Value = 0.6061951661742485
$ws.Range("D26").Value = 0.6290105436377385
$ws.Range("E26").Value = 0.579735806092556
$ws.Range("F26").Value = 0.6295883380069857
$ws.Range("G26").Value = 0.6466451715177792
$ws.Range("H26").Value = 0.6422907871590777
$ws.Range("I26").Value = 0.5460161349399134
$ws.Range("K26").Value = 0.6139741398787006
$ws.Range("L26").Value = 0.6246909267948835
$ws.Range("M26").Value = 0.6666421500649489
$ws.Range("N26").Value = 0.6279718323457103
$ws.Range("O26").Value = 0.6827780440518729
$ws.Range("P26").Value = 0.6496462751434467
$ws.Range("Q26").Value = 0.5349718102811883
$ws.Range("R26").Value = 0.6694420657237685
$ws.Range("K27").Value = 0.6241150885277071
$ws.Range("C29").Value = 0.6443822364380736
$ws.Range("D29").Value = 0.6894203090619051
$ws.Range("E29").Value = 0.6401029727022728
$ws.Range("F29").Value = 0.6868823760541658
$ws.Range("G29").Value = 0.6869629707823841
$ws.Range("H29").Value = 0.6087381451283874
$ws.Range("I29").Value = 0.548786621744655
$ws.Range("K29").Value = 0.6504904013681438
$ws.Range("L29").Value = 0.6226099391489784
$ws.Range("M29").Value = 0.6059129560318404
$ws.Range("N29").Value = 0.558612612254482
$ws.Range("O29").Value = 0.5565778251872112
$ws.Range("P29").Value = 0.5975144050157137
$ws.Range("Q29").Value = 0.6189178305608056
$ws.Range("R29").Value = 0.7306215018151387
$ws.Range("C34").Value = 0.8234993659411382
$ws.Range("D34").Value = 0.6519068856037634
$ws.Range("E34").Value = 0.6240810970156861
$ws.Range("F34").Value = 0.7195103674033672
$ws.Range("C35").Value = 0.6085330998540964
$ws.Range("D35").Value = 0.6078074062845437
$ws.Range("E35").Value = 0.5898409602227346
$ws.Range("F35").Value = 0.6345678365882789
$ws.Range("C36").Value = 0.6182350050858616
$ws.Range("D36").Value = 0.600760353992564
$ws.Range("E36").Value = 0.6503458456801725
$ws.Range("F36").Value = 0.6022069380024784
$ws.Range("D37").Value = 0.5502122451073774
$ws.Range("C39").Value = 0.6182337753633009
$ws.Range("D39").Value = 0.6009598282345683
$ws.Range("E39").Value = 0.5675454481434694
$ws.Range("F39").Value = 0.6520950473313748

# Row 24: move the bold/underlined "max" highlight from D24 to F24
$ws.Range("D24").Font.Bold = $false
$ws.Range("D24").Font.Underline = $false
$ws.Range("F24").Font.Bold = $true
$ws.Range("F24").Font.Underline = $true
